# "fresh start with ParticleCube"
# Wipe out all the old notes/content from the document, leaving only the
# single trailing empty paragraph that precedes the section properties.
$d = $word.ActiveDocument

while ($d.Paragraphs.Count -gt 1) {
    $p = $d.Paragraphs.Item(1)
    $p.Range.Delete()
}
